$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7000.5
$ws.Range("I18").Value = 4001
$ws.Range("K18").Value = 4001
$ws.Range("M18").Value = -3717
$ws.Range("H19").Value = 1371.5333
$ws.Range("I19").Value = 1281.2727
$ws.Range("K19").Value = 1281.2727
$ws.Range("M19").Value = -1106.2727
$ws.Range("H106").Value = 4387.7144
$ws.Range("I106").Value = 4519
$ws.Range("J106").Value = 3600
$ws.Range("K106").Value = 4519
$ws.Range("L106").Value = 3600
$ws.Range("M106").Value = -3888
$ws.Range("N106").Value = -4862
$ws.Range("H116").Value = 3984.5833
$ws.Range("I116").Value = 3983.7
$ws.Range("K116").Value = 3983.7
$ws.Range("M116").Value = -541.6999999999998
$ws.Range("H132").Value = 1151.2778
$ws.Range("I132").Value = 1218.1177
$ws.Range("K132").Value = 3654.3531
$ws.Range("M132").Value = -1124.3531
$ws.Range("H135").Value = 2686.9092
$ws.Range("I135").Value = 3177.625
$ws.Range("J135").Value = 1378.3334
$ws.Range("K135").Value = 28598.625
$ws.Range("L135").Value = 12405.0006
$ws.Range("M135").Value = -26063.625
$ws.Range("N135").Value = -17475.0006
$ws.Range("H137").Value = 2436.3333
$ws.Range("I137").Value = 1567.875
$ws.Range("K137").Value = 4703.625
$ws.Range("M137").Value = -2153.625
$ws.Range("H138").Value = 6892.567
$ws.Range("J138").Value = 7846.625
$ws.Range("L138").Value = 23539.875
$ws.Range("N138").Value = -33819.875
$ws.Range("H141").Value = 5848.25
$ws.Range("I141").Value = 4969.4287
$ws.Range("J141").Value = 12000
$ws.Range("K141").Value = 14908.2861
$ws.Range("L141").Value = 36000
$ws.Range("M141").Value = -9728.286100000001
$ws.Range("N141").Value = -46360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 590
$ws.Range("I19").Value = 590
$ws.Range("K19").Value = 590
$ws.Range("M19").Value = -361
$ws.Range("H32").Value = 4326.1816
$ws.Range("I32").Value = 3395.6453
$ws.Range("J32").Value = 18749.5
$ws.Range("K32").Value = 3395.6453
$ws.Range("L32").Value = 18749.5
$ws.Range("M32").Value = -3108.6453
$ws.Range("N32").Value = -19323.5
$ws.Range("H132").Value = 3046.9375
$ws.Range("I132").Value = 2657.75
$ws.Range("K132").Value = 7973.25
$ws.Range("M132").Value = -5443.25
$ws.Range("H133").Value = 72500
$ws.Range("J133").Value = 72500
$ws.Range("L133").Value = 72500
$ws.Range("N133").Value = -77560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4907.2188
$ws.Range("J31").Value = 6611.0625
$ws.Range("L31").Value = 6611.0625
$ws.Range("N31").Value = -7201.0625
$ws.Range("H34").Value = 4907.2188
$ws.Range("J34").Value = 6611.0625
$ws.Range("L34").Value = 6611.0625
$ws.Range("N34").Value = -7015.0625
$ws.Range("H58").Value = 2334.75
$ws.Range("J58").Value = 2387
$ws.Range("L58").Value = 2387
$ws.Range("N58").Value = -2793
$ws.Range("H99").Value = 1508.5
$ws.Range("I99").Value = 1633
$ws.Range("K99").Value = 1633
$ws.Range("M99").Value = -135
$ws.Range("H126").Value = 1508.5
$ws.Range("I126").Value = 1633
$ws.Range("K126").Value = 4899
$ws.Range("M126").Value = -2429
$ws.Range("H132").Value = 4346.6
$ws.Range("I132").Value = 3883.6155
$ws.Range("J132").Value = 5206.4287
$ws.Range("K132").Value = 11650.8465
$ws.Range("L132").Value = 15619.2861
$ws.Range("M132").Value = -9120.8465
$ws.Range("N132").Value = -20679.2861
$ws.Range("H134").Value = 8461.625
$ws.Range("I134").Value = 8461.625
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 25384.875
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -22849.875
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2334.75
$ws.Range("J136").Value = 2387
$ws.Range("L136").Value = 7161
$ws.Range("N136").Value = -12261

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3086461.8
$ws.Range("J2").Value = 80
$ws.Range("L2").Value = 480
$ws.Range("N2").Value = -706
$ws.Range("H14").Value = 345.5
$ws.Range("I14").Value = 345.5
$ws.Range("K14").Value = 1036.5
$ws.Range("M14").Value = -863.5
$ws.Range("H34").Value = 2448.3333
$ws.Range("J34").Value = 3463.3333
$ws.Range("L34").Value = 10389.9999
$ws.Range("N34").Value = -10557.9999
$ws.Range("H68").Value = 1267.3334
$ws.Range("J68").Value = 1301
$ws.Range("L68").Value = 3903
$ws.Range("N68").Value = -5525
$ws.Range("H71").Value = 1267.3334
$ws.Range("J71").Value = 1301
$ws.Range("L71").Value = 11709
$ws.Range("N71").Value = -19821
$ws.Range("H113").Value = 624.3
$ws.Range("J113").Value = 787
$ws.Range("L113").Value = 2361
$ws.Range("N113").Value = -6701
$ws.Range("H117").Value = 441.33334
$ws.Range("I117").Value = 212.5
$ws.Range("K117").Value = 637.5
$ws.Range("M117").Value = 2804.5
$ws.Range("H122").Value = 3819.0454
$ws.Range("J122").Value = 3945.2144
$ws.Range("L122").Value = 35506.9296
$ws.Range("N122").Value = -40406.9296
$ws.Range("H129").Value = 432
$ws.Range("I129").Value = 432
$ws.Range("K129").Value = 1296
$ws.Range("M129").Value = 3704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 31.857143
$ws.Range("I2").Value = 20.5
$ws.Range("K2").Value = 20.5
$ws.Range("M2").Value = 92.5
$ws.Range("H80").Value = 3417.5217
$ws.Range("I80").Value = 2612.2856
$ws.Range("K80").Value = 2612.2856
$ws.Range("M80").Value = -1614.2856
$ws.Range("H83").Value = 3417.5217
$ws.Range("I83").Value = 2612.2856
$ws.Range("K83").Value = 13061.428
$ws.Range("M83").Value = -8069.428
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4422.421
$ws.Range("I132").Value = 4285.7334
$ws.Range("K132").Value = 12857.2002
$ws.Range("M132").Value = -10327.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4999.5
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4798
$ws.Range("H82").Value = 911.875
$ws.Range("I82").Value = 911.875
$ws.Range("K82").Value = 911.875
$ws.Range("M82").Value = -550.875
$ws.Range("H85").Value = 911.875
$ws.Range("I85").Value = 911.875
$ws.Range("K85").Value = 911.875
$ws.Range("M85").Value = 336.125
$ws.Range("H113").Value = 4999.5
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
$ws.Range("H132").Value = 5598.2
$ws.Range("J132").Value = 5746.75
$ws.Range("L132").Value = 17240.25
$ws.Range("N132").Value = -22300.25
$ws.Range("H136").Value = 32905.625
$ws.Range("I136").Value = 4833.3335
$ws.Range("K136").Value = 14500.0005
$ws.Range("M136").Value = -11950.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 800.6
$ws.Range("I113").Value = 738.1429000000001
$ws.Range("J113").Value = 946.3333
$ws.Range("K113").Value = 2214.4287
$ws.Range("L113").Value = 2838.9999
$ws.Range("M113").Value = -44.42870000000039
$ws.Range("N113").Value = -7178.9999
$ws.Range("H132").Value = 3491.4736
$ws.Range("I132").Value = 3189.8667
$ws.Range("K132").Value = 9569.6001
$ws.Range("M132").Value = -7039.6001
